$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.231.12'
$ws.Range("E2").Value = '  +2.41%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.425.28'
$ws.Range("E3").Value = '  +1.94%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.42'
$ws.Range("E5").Value = '  +1.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.65'
$ws.Range("E6").Value = '  +3.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.513'
$ws.Range("E7").Value = '  +0.80%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.499'
$ws.Range("E9").Value = '  -0.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.32'
$ws.Range("E10").Value = '  +2.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0801'
$ws.Range("E11").Value = '  +1.65%  '
$ws.Range("E12").Value = '  +2.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.81'
$ws.Range("E13").Value = '  +2.31%  '
$ws.Range("E14").Value = '  +1.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.803.78'
$ws.Range("E15").Value = '  +1.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.433.94'
$ws.Range("E16").Value = '  +2.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.835'
$ws.Range("E17").Value = '  +3.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '44.181.44'
$ws.Range("E18").Value = '  +2.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.30'
$ws.Range("E19").Value = '  +0.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.41'
$ws.Range("E20").Value = '  +1.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0905'
$ws.Range("E21").Value = '  +1.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.46'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.30'
$ws.Range("E23").Value = '  +3.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '240.46'
$ws.Range("E24").Value = '  +2.05%  '
$ws.Range("E25").Value = '  +1.13%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.23'
$ws.Range("E27").Value = '  +1.52%  '
$ws.Range("E28").Value = '  -0.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.59'
$ws.Range("E29").Value = '  +4.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.08'
$ws.Range("E30").Value = '  +4.25%  '
$ws.Range("B31").Value = 'Celestia'
$ws.Range("C31").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '18.63'
$ws.Range("E31").Value = '  +7.72%  '
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.116'
$ws.Range("E32").Value = '  +10.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.17'
$ws.Range("E33").Value = '  +1.43%  '
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0759'
$ws.Range("E35").Value = '  +2.41%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.90'
$ws.Range("E36").Value = '  +2.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.48'
$ws.Range("E37").Value = '  +4.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '128.45'
$ws.Range("E38").Value = '  +23.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.90'
$ws.Range("E39").Value = '  +3.69%  '
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.10'
$ws.Range("E42").Value = '  -6.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0287'
$ws.Range("E43").Value = '  +2.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.954.80'
$ws.Range("E44").Value = '  -0.21%  '
$ws.Range("E45").Value = '  +1.81%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.87'
$ws.Range("E46").Value = '  +4.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.41'
$ws.Range("E47").Value = '  +2.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.38'
$ws.Range("E49").Value = '  +1.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '73.61'
$ws.Range("E50").Value = '  +2.29%  '
$ws.Range("E51").Value = '  +0.96%  '
